$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.62
$ws.Range("I2").Value = 5.5
$ws.Range("Q2").Value = 1.91
$ws.Range("R2").Value = 1.99
$ws.Range("U2").Value = 1.87
$ws.Range("V2").Value = 1.87
$ws.Range("AH2").Value = 13
$ws.Range("AI2").Value = 26
$ws.Range("AO2").Value = 8.5
$ws.Range("O3").Value = 1.62
$ws.Range("P3").Value = 2.3
$ws.Range("Q3").Value = 2.88
$ws.Range("R3").Value = 1.4
$ws.Range("I4").Value = 6
$ws.Range("G5").Value = 3
$ws.Range("I5").Value = 2.3
$ws.Range("J5").Value = 3.6
$ws.Range("Z5").Value = 34
$ws.Range("AB5").Value = 34
$ws.Range("AK5").Value = 21
$ws.Range("AR5").Value = 81
$ws.Range("AS5").Value = 201
$ws.Range("AW5").Value = 4.33
$ws.Range("M7").Value = 1.02
$ws.Range("N7").Value = 7.4
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 3.9
$ws.Range("J8").Value = 2.63
$ws.Range("L8").Value = 4
$ws.Range("N8").Value = 10
$ws.Range("U8").Value = 1.73
$ws.Range("V8").Value = 2
$ws.Range("X8").Value = 9.5
$ws.Range("Y8").Value = 9
$ws.Range("Z8").Value = 17
$ws.Range("AA8").Value = 17
$ws.Range("AD8").Value = 6
$ws.Range("AE8").Value = 13
$ws.Range("AF8").Value = 41
$ws.Range("AG8").Value = 201
$ws.Range("AJ8").Value = 13
$ws.Range("AL8").Value = 29
$ws.Range("AO8").Value = 11
$ws.Range("AQ8").Value = 41
$ws.Range("AW8").Value = 5.5
$ws.Range("BA8").Value = 81
$ws.Range("M9").Value = 1.03
$ws.Range("O9").Value = 1.14
$ws.Range("Q9").Value = 1.53
$ws.Range("R9").Value = 2.4
$ws.Range("M10").Value = 1.08
$ws.Range("O10").Value = 1.4
$ws.Range("Q10").Value = 2.3
$ws.Range("R10").Value = 1.6
$ws.Range("M11").Value = 1.13
$ws.Range("O11").Value = 1.57
$ws.Range("W11").Value = 5
$ws.Range("AE11").Value = 23
$ws.Range("AF11").Value = 101
$ws.Range("H13").Value = 5.25
$ws.Range("K13").Value = 2.63
$ws.Range("M13").Value = 1.04
$ws.Range("N13").Value = 13
$ws.Range("O13").Value = 1.18
$ws.Range("P13").Value = 4.5
$ws.Range("Q13").Value = 1.62
$ws.Range("R13").Value = 2.25
$ws.Range("W13").Value = 7.5
$ws.Range("Y13").Value = 9.5
$ws.Range("AA13").Value = 11
$ws.Range("AC13").Value = 13
$ws.Range("AO13").Value = 5.5
